$wb = $excel.ActiveWorkbook

# Clear the old Phone/Price header row on Sheet1
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1:B1").ClearContents()

# Add a new worksheet named "data" at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "data"
$newSheet.Range("A2").Value = "Apple iPhone X (Silver, 64 GB)"

# Keep Sheet1 as the active/selected sheet
$ws1.Activate()
